$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell G1, copying the style from F1 (bold header style)
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "G"

# Fill G2:G105 with data values
$ws.Range("G2").Value = "Hb 47"
$ws.Range("G3").Value = "Hb 48"
$ws.Range("G4").Value = "Hb 21"
$ws.Range("G5").Value = "Hb 22"
$ws.Range("G6").Value = "Hb 17"
$ws.Range("G7").Value = "S 6"
$ws.Range("G8").Value = "Hb 7"
$ws.Range("G9").Value = "Hb 46"
$ws.Range("G10").Value = "Hb 1"
$ws.Range("G11").Value = "Hb 2"
$ws.Range("G12").Value = "Hb 3"
$ws.Range("G13").Value = "Hb 5"
$ws.Range("G14").Value = "S 24"
$ws.Range("G15").Value = "S 25"
$ws.Range("G16").Value = "S 26"
$ws.Range("G17").Value = "S 27"
$ws.Range("G18").Value = "S 28"
$ws.Range("G19").Value = "Hb 103"
$ws.Range("G20").Value = "Hb 104"
$ws.Range("G21").Value = "Hb 105"
$ws.Range("G22").Value = "Hb 106"
$ws.Range("G23").Value = "Hb 107"
$ws.Range("G24").Value = "Hb 63"
$ws.Range("G25").Value = "Hb 65"
$ws.Range("G26").Value = "Hb 66"
$ws.Range("G27").Value = "Hb 67"
$ws.Range("G28").Value = "Hb 68"
$ws.Range("G29").Value = "Hb 69"
$ws.Range("G30").Value = "Hb 70"
$ws.Range("G31").Value = "Hb 71"
$ws.Range("G32").Value = "Hb 72"
$ws.Range("G33").Value = "Hb 93"
$ws.Range("G34").Value = "Hb 94"
$ws.Range("G35").Value = "Hb 95"
$ws.Range("G36").Value = "Hb 96"
$ws.Range("G37").Value = "Hb 97"
$ws.Range("G38").Value = "Hb 98"
$ws.Range("G39").Value = "Hb 99"
$ws.Range("G40").Value = "Hb 100"
$ws.Range("G41").Value = "Hb 101"
$ws.Range("G42").Value = "Hb 102"
$ws.Range("G43").Value = "S 29"
$ws.Range("G44").Value = "S 30"
$ws.Range("G45").Value = "Hb 83"
$ws.Range("G46").Value = "Hb 84"
$ws.Range("G47").Value = "Hb 85"
$ws.Range("G48").Value = "Hb 86"
$ws.Range("G49").Value = "Hb 87"
$ws.Range("G50").Value = "Hb 88"
$ws.Range("G51").Value = "Hb 89"
$ws.Range("G52").Value = "Hb 90"
$ws.Range("G53").Value = "Hb 91"
$ws.Range("G54").Value = "Hb 92"
$ws.Range("G55").Value = "Hb 40"
$ws.Range("G56").Value = "Hb 41"
$ws.Range("G57").Value = "Hb 42"
$ws.Range("G58").Value = "Hb 43"
$ws.Range("G59").Value = "S 8"
$ws.Range("G60").Value = "S 9"
$ws.Range("G61").Value = "S 11"
$ws.Range("G62").Value = "S 12"
$ws.Range("G63").Value = "Hb 53"
$ws.Range("G64").Value = "Hb 54"
$ws.Range("G65").Value = "Hb 55"
$ws.Range("G66").Value = "Hb 56"
$ws.Range("G67").Value = "Hb 57"
$ws.Range("G68").Value = "Hb 58"
$ws.Range("G69").Value = "Hb 59"
$ws.Range("G70").Value = "Hb 60"
$ws.Range("G71").Value = "Hb 61"
$ws.Range("G72").Value = "Hb 62"
$ws.Range("G73").Value = "Hb 12"
$ws.Range("G74").Value = "Hb 13"
$ws.Range("G75").Value = "Hb 14"
$ws.Range("G76").Value = "Hb 15"
$ws.Range("G77").Value = "S 18"
$ws.Range("G78").Value = "S 19"
$ws.Range("G79").Value = "S 21"
$ws.Range("G80").Value = "S 22"
$ws.Range("G81").Value = "Hb 35"
$ws.Range("G82").Value = "Hb 36"
$ws.Range("G83").Value = "Hb 38"
$ws.Range("G84").Value = "Hb 39"
$ws.Range("G85").Value = "S 1"
$ws.Range("G86").Value = "S 2"
$ws.Range("G87").Value = "S 3"
$ws.Range("G88").Value = "S 4"
$ws.Range("G89").Value = "S 5"
$ws.Range("G90").Value = "Hb 73"
$ws.Range("G91").Value = "Hb 74"
$ws.Range("G92").Value = "Hb 75"
$ws.Range("G93").Value = "Hb 76"
$ws.Range("G94").Value = "Hb 77"
$ws.Range("G95").Value = "Hb 78"
$ws.Range("G96").Value = "Hb 79"
$ws.Range("G97").Value = "Hb 31"
$ws.Range("G98").Value = "Hb 32"
$ws.Range("G99").Value = "Hb 33"
$ws.Range("G100").Value = "S 13"
$ws.Range("G101").Value = "S 14"
$ws.Range("G102").Value = "S 15"
$ws.Range("G103").Value = "S 16"
$ws.Range("G104").Value = "Hb 10"
$ws.Range("G105").Value = "Hb 11"
